$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'27.754.59"
$c.Style = "Normal"

$c = $ws.Range("E2")
$c.Value = "'  +0.60%  "
$c.Style = "Normal"

$c = $ws.Range("D3")
$c.Value = "'1.638.34"
$c.Style = "Normal"

$c = $ws.Range("E4")
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.Value = "'212.57"
$c.Style = "Normal"

$c = $ws.Range("E6")
$c.Value = "'  -1.91%  "
$c.Style = "Normal"

$c = $ws.Range("E7")
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D8")
$c.Value = "'23.25"
$c.Style = "Normal"

$c = $ws.Range("E8")
$c.Value = "'  -1.16%  "
$c.Style = "Normal"

$c = $ws.Range("E9")
$c.Value = "'  +2.16%  "
$c.Style = "Normal"

$c = $ws.Range("E10")
$c.Value = "'  -0.05%  "
$c.Style = "Normal"

$c = $ws.Range("D11")
$c.Value = "'0.0889"
$c.Style = "Normal"

$c = $ws.Range("E11")
$c.Value = "'  +0.14%  "
$c.Style = "Normal"

$c = $ws.Range("D12")
$c.Value = "'1.871.31"
$c.Style = "Normal"

$c = $ws.Range("E12")
$c.Value = "'  -0.51%  "
$c.Style = "Normal"

$c = $ws.Range("D13")
$c.Value = "'1.639.79"
$c.Style = "Normal"

$c = $ws.Range("E13")
$c.Value = "'  -0.51%  "
$c.Style = "Normal"

$c = $ws.Range("E14")
$c.Value = "'  +0.30%  "
$c.Style = "Normal"

$c = $ws.Range("D15")
$c.Value = "'0.562"
$c.Style = "Normal"

$c = $ws.Range("E15")
$c.Value = "'  -3.87%  "
$c.Style = "Normal"

$c = $ws.Range("D16")
$c.Value = "'64.76"
$c.Style = "Normal"

$c = $ws.Range("E16")
$c.Value = "'  +0.52%  "
$c.Style = "Normal"

$c = $ws.Range("D17")
$c.Value = "'27.723.66"
$c.Style = "Normal"

$c = $ws.Range("E17")
$c.Value = "'  +0.61%  "
$c.Style = "Normal"

$c = $ws.Range("D18")
$c.Value = "'229.89"
$c.Style = "Normal"

$c = $ws.Range("E18")
$c.Value = "'  -0.64%  "
$c.Style = "Normal"

$c = $ws.Range("D19")
$c.Value = "'7.70"
$c.Style = "Normal"

$c = $ws.Range("E19")
$c.Value = "'  +1.83%  "
$c.Style = "Normal"

$c = $ws.Range("D20")
$c.Value = "'0.0₃0722"
$c.Style = "Normal"

$c = $ws.Range("E20")
$c.Value = "'  -0.04%  "
$c.Style = "Normal"

$c = $ws.Range("E21")
$c.Value = "'  -0.09%  "
$c.Style = "Normal"

$c = $ws.Range("D22")
$c.Value = "'4.31"
$c.Style = "Normal"

$c = $ws.Range("E22")
$c.Value = "'  -0.33%  "
$c.Style = "Normal"

$c = $ws.Range("D23")
$c.Value = "'10.20"
$c.Style = "Normal"

$c = $ws.Range("E23")
$c.Value = "'  +4.95%  "
$c.Style = "Normal"

$c = $ws.Range("D24")
$c.Value = "'2.11"
$c.Style = "Normal"

$c = $ws.Range("E24")
$c.Value = "'  +4.96%  "
$c.Style = "Normal"

$c = $ws.Range("D25")
$c.Value = "'151.20"
$c.Style = "Normal"

$c = $ws.Range("E25")
$c.Value = "'  +1.54%  "
$c.Style = "Normal"

$c = $ws.Range("D26")
$c.Value = "'6.94"
$c.Style = "Normal"

$c = $ws.Range("E26")
$c.Value = "'  -0.86%  "
$c.Style = "Normal"

$c = $ws.Range("E27")
$c.Value = "'  -0.63%  "
$c.Style = "Normal"

$c = $ws.Range("E28")
$c.Value = "'  -0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D29")
$c.Value = "'15.60"
$c.Style = "Normal"

$c = $ws.Range("E29")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("E30")
$c.Value = "'  +0.42%  "
$c.Style = "Normal"

$c = $ws.Range("D31")
$c.Value = "'0.0486"
$c.Style = "Normal"

$c = $ws.Range("E31")
$c.Value = "'  -0.06%  "
$c.Style = "Normal"

$c = $ws.Range("E32")
$c.Value = "'  +0.07%  "
$c.Style = "Normal"

$c = $ws.Range("D33")
$c.Value = "'1.471.07"
$c.Style = "Normal"

$c = $ws.Range("E33")
$c.Value = "'  +3.16%  "
$c.Style = "Normal"

$c = $ws.Range("D34")
$c.Value = "'3.11"
$c.Style = "Normal"

$c = $ws.Range("E34")
$c.Value = "'  -1.83%  "
$c.Style = "Normal"

$c = $ws.Range("E35")
$c.Value = "'  -2.47%  "
$c.Style = "Normal"

$c = $ws.Range("E36")
$c.Value = "'  -0.53%  "
$c.Style = "Normal"

$c = $ws.Range("E37")
$c.Value = "'  -0.52%  "
$c.Style = "Normal"

$c = $ws.Range("D38")
$c.Value = "'0.881"
$c.Style = "Normal"

$c = $ws.Range("E38")
$c.Value = "'  -0.16%  "
$c.Style = "Normal"

$c = $ws.Range("E39")
$c.Value = "'  +0.11%  "
$c.Style = "Normal"

$c = $ws.Range("D40")
$c.Value = "'0.913"
$c.Style = "Normal"

$c = $ws.Range("E40")
$c.Value = "'  +12.43%  "
$c.Style = "Normal"

$c = $ws.Range("D41")
$c.Value = "'68.84"
$c.Style = "Normal"

$c = $ws.Range("E41")
$c.Value = "'  +5.76%  "
$c.Style = "Normal"

$c = $ws.Range("E43")
$c.Value = "'  -1.82%  "
$c.Style = "Normal"

$c = $ws.Range("D44")
$c.Value = "'5.60"
$c.Style = "Normal"

$c = $ws.Range("E44")
$c.Value = "'  +1.34%  "
$c.Style = "Normal"

$c = $ws.Range("E45")
$c.Value = "'  -0.94%  "
$c.Style = "Normal"

$c = $ws.Range("E46")
$c.Value = "'  -0.59%  "
$c.Style = "Normal"

$c = $ws.Range("D47")
$c.Value = "'1.780.85"
$c.Style = "Normal"

$c = $ws.Range("E47")
$c.Value = "'  -0.54%  "
$c.Style = "Normal"

$c = $ws.Range("D48")
$c.Value = "'1.71"
$c.Style = "Normal"

$c = $ws.Range("E48")
$c.Value = "'  +1.86%  "
$c.Style = "Normal"

$c = $ws.Range("D49")
$c.Value = "'86.93"
$c.Style = "Normal"

$c = $ws.Range("E49")
$c.Value = "'  -1.46%  "
$c.Style = "Normal"

$c = $ws.Range("E50")
$c.Value = "'  -1.06%  "
$c.Style = "Normal"

$c = $ws.Range("E51")
$c.Value = "'  -0.19%  "
$c.Style = "Normal"

